# Update vm_pu results sheet: case with 380 kV done
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @"
2|1.02,1.061580654564151,1.065827896789079,1.057131682040233,1.069972722833207|1.027727093623883,1.066555599747767,1.068540087318213,1.059867491236546,1.072673802253232,1.068070231035047
3|1.02,1.065040037956285,1.069068234466548,1.060135155210381,1.073408141163539|1.028070710856076,1.069657140357252,1.071589856714869,1.062679185340305,1.07591902199196,1.071176176187948
4|1.02,1.067248008602524,1.071134168909404,1.062052334200652,1.075592677910512|1.028281227693845,1.07163410748539,1.073532413656401,1.064472152701901,1.077980474568658,1.073155950836169
5|1.02,1.068169144848626,1.071995513156468,1.062852201119229,1.076502083246614|1.02836692821246,1.072458245337372,1.07434186770134,1.065219767329366,1.078838124967922,1.073981259058431
6|1.02,1.068323397383298,1.072139721420079,1.062986148615462,1.076654256415823|1.02838115432016,1.072596217761745,1.074477361885809,1.065344939300488,1.078981607441688,1.074119427419472
7|1.02,1.067260344483711,1.071145706172244,1.062063045864166,1.075604864403499|1.028282383793984,1.071645146834968,1.073543257655601,1.064482166307488,1.077991969559775,1.07316700586289
8|1.02,1.062756205664002,1.066929477155361,1.058152265983583,1.071141827687362|1.027845688427606,1.067610094292387,1.06957727431475,1.060823283544782,1.073778626852897,1.069126223082872
9|1.02,1.05457633822253,1.059255038775264,1.051051604810684,1.062972773626579|1.026984124997332,1.06026168115583,1.062343576198544,1.054165851622549,1.066049802426094,1.061767374356786
10|1.02,1.048945691897045,1.053960614735169,1.046165049848122,1.057306531186098|1.026345641756164,1.055189612003984,1.057343246050095,1.049574718749255,1.06067769226162,1.05668810228497
11|1.02,1.046462152304471,1.05162261630905,1.044010026579713,1.054797083036346|1.026053424529779,1.052949190919478,1.055132754292324,1.047547687439034,1.058295873417507,1.054444499545544
12|1.02,1.045532528053492,1.05074705687994,1.043203420550826,1.05385622645859|1.025942465822295,1.052110080614742,1.054304589554308,1.046788641848603,1.05740247446871,1.053604197607928
13|1.02,1.045732263146218,1.050935194503871,1.043376722179731,1.054058444261204|1.025966377019747,1.052290390357483,1.0544825592172,1.04695174063664,1.057594509998905,1.053784763411185
14|1.02,1.046385456464021,1.051550389328819,1.043943478956713,1.054719491739061|1.026044302239201,1.052879972603673,1.055064444246472,1.047485070753219,1.058222203893707,1.054375182931792
15|1.02,1.046786956730169,1.051928478387151,1.04429185568088,1.055125617043827|1.026091992815471,1.053242308315352,1.055422014684583,1.047812854898847,1.058607785948017,1.054738033201742
16|1.02,1.04910953367824,1.054114797416735,1.046307225864391,1.057471868637501|1.026364699228755,1.055337346933546,1.057488970681316,1.049708402530066,1.060834565879851,1.056836047015077
17|1.02,1.05055405360241,1.055473841964218,1.047560764834458,1.058928407358218|1.026531509826261,1.056639487478737,1.058773191354055,1.050886804631241,1.062216241973179,1.0581400367492
18|1.02,1.051392255249807,1.056262182124983,1.048288177261129,1.059772609015885|1.026627291553075,1.057394760871275,1.05951790304564,1.051570397428649,1.063016802684503,1.058896382716586
19|1.02,1.051677329297004,1.056530253489189,1.04853557668875,1.060059558470135|1.02665969495871,1.057651578128976,1.059771100521481,1.051802856359209,1.063288875670698,1.059153564684251
20|1.02,1.050399523633523,1.055328483147265,1.047426662252292,1.058772692614951|1.026513769793371,1.056500220988964,1.058635858814459,1.050760762804018,1.062068556607348,1.058000572485037
21|1.02,1.046193306557141,1.051369428819042,1.043776754696832,1.054525073655253|1.026021422321847,1.052706548631622,1.054893291657533,1.047328189236921,1.058037606133011,1.054201512677808
22|1.02,1.043507294596701,1.048838848916625,1.04144628278001,1.051803705481515|1.025697857253487,1.050281139700268,1.052499021570306,1.045134475064654,1.055452763139123,1.051772659387574
23|1.02,1.044935229208969,1.050184379272343,1.04268517729243,1.0532512783055|1.025870730449888,1.051570800829639,1.053772269954542,1.046300858512444,1.056827929706489,1.053064151983674
24|1.02,1.050469362534804,1.055394178016388,1.047487269024836,1.058843069984619|1.02652179043213,1.056563162613996,1.058697926908439,1.050817727248967,1.062135305769897,1.058063603494396
25|1.02,1.056721224928137,1.061269430093215,1.052913305001949,1.065122301412125|1.027217986975638,1.062190939807917,1.064244010205935,1.05591300574543,1.068085443372787,1.063699372777448
"@

$lines = $data -split "`n"
foreach ($line in $lines) {
    $line = $line.Trim()
    if ($line -eq "") { continue }
    $parts = $line -split '\|'
    $rowNum = [int]$parts[0]
    $bfVals = $parts[1] -split ','
    $inVals = $parts[2] -split ','

    $bfArr = New-Object 'object[,]' 1,$bfVals.Length
    for ($i = 0; $i -lt $bfVals.Length; $i++) {
        $bfArr[0,$i] = [double]$bfVals[$i]
    }
    $ws.Range("B$rowNum`:F$rowNum").Value = $bfArr

    $inArr = New-Object 'object[,]' 1,$inVals.Length
    for ($i = 0; $i -lt $inVals.Length; $i++) {
        $inArr[0,$i] = [double]$inVals[$i]
    }
    $ws.Range("I$rowNum`:N$rowNum").Value = $inArr
}
